$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("C4").Value = 1555.2
$ws1.Range("M29").Value = 1672.48
$ws1.Range("C55").Value = "2 de 53"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 4324.38
$ws2.Range("F29").Value = 1672.48
$ws2.Range("F55").Value = 23839.89

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# The COM layer adds ~5/6 (one default char) of padding when reading back
# ColumnWidth vs. the raw OOXML <col width>, so subtract it to land on 24.
$ws3.Columns.Item(6).ColumnWidth = 24 - (5/6)

$ws3.Range("D2").Value = 2073.6
$ws3.Range("E2").Value = 7896.74304517915
$ws3.Range("F2").Value = 0.2079767958438125

$ws3.Range("D16").Value = 21005.76
$ws3.Range("E16").Value = 30820.7
$ws3.Range("F16").Value = 0.4053095658086622

$ws3.Range("D19").Value = 23839.89
$ws3.Range("E19").Value = 89866.56064517915
$ws3.Range("F19").Value = 0.2096617198472965
